# Issue #57 - Make genre required with PBCore controlled vocabulary.
# Fixture update: add a new "Genre" override column (U) to the batch-ingest
# manifest spreadsheet, modeled on the existing duplicate "Topical Subject"
# column (T) used to let the manifest spreadsheet override file-level
# metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header, mirroring the existing "Genre" header used in D2.
$ws.Range("U2").Value = "Genre"

# New per-file Genre override values (PBCore controlled-vocabulary terms)
# for the two file rows already present in the fixture.
$ws.Range("U3").Value = "Auction"
$ws.Range("U4").Value = "Anime"

# Move the active selection to the newly added cell, matching the
# spreadsheet author's last edit position.
$ws.Range("U4").Select()
